$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the SmartScore cells as plain text (matching the Streamlit export,
#     which writes these as text rather than numbers) before touching any other cell,
#     so that later row auto-height recalculation does not leave stray style indices.
$ws.Range('I11').NumberFormat = '@'
$ws.Range('I11').Value = '0.579'
$ws.Range('I11').ClearFormats()
$ws.Range('L11').NumberFormat = '@'
$ws.Range('L11').Value = '0.479'
$ws.Range('L11').ClearFormats()
$ws.Range('O11').NumberFormat = '@'
$ws.Range('O11').Value = '0.469'
$ws.Range('O11').ClearFormats()
$ws.Range('R11').NumberFormat = '@'
$ws.Range('R11').Value = '0.601'
$ws.Range('R11').ClearFormats()
$ws.Range('U11').NumberFormat = '@'
$ws.Range('U11').Value = '0.559'
$ws.Range('U11').ClearFormats()
$ws.Range('X11').NumberFormat = '@'
$ws.Range('X11').Value = '0.547'
$ws.Range('X11').ClearFormats()
$ws.Range('AA11').NumberFormat = '@'
$ws.Range('AA11').Value = '0.715'
$ws.Range('AA11').ClearFormats()
$ws.Range('AD11').NumberFormat = '@'
$ws.Range('AD11').Value = '0.590'
$ws.Range('AD11').ClearFormats()
$ws.Range('AG11').NumberFormat = '@'
$ws.Range('AG11').Value = '0.567'
$ws.Range('AG11').ClearFormats()

# --- Step 2: write the remaining cells for the new row (row 11)
$ws.Range('A11').Value = 'Alys_20251120_203615'
$ws.Range('C11').Value = 'Alys'
$ws.Range("D11").Value = 20
$ws.Range('E11').Value = 'Female'
$ws.Range('F11').Value = '2025-11-20 20:36:15'
$ws.Range('G11').Value = '{
  "portion": 0.8,
  "diet": 0.8571428571428571,
  "salt": 0.4,
  "fat": 1.0,
  "natural": 0.8,
  "convenience": 0.8,
  "price": 0.8
}'
$ws.Range('H11').Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range('J11').Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range('K11').Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range('M11').Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range('N11').Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range('P11').Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range('Q11').Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range('S11').Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range('T11').Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range('V11').Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range('W11').Value = 'Annie’s Shells & White Cheddar'
$ws.Range('Y11').Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range('Z11').Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range('AB11').Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range('AC11').Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range('AE11').Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range('AF11').Value = 'Jack Link’s Beef Jerky Original'
$ws.Range('AH11').Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

# Excel recalculates a larger row height for wrapped/multi-line content; AutoFit
# restores the row to the sheet's default (non-custom) height, matching the source.
$ws.Rows("11").AutoFit()
